# Leakage test / Jankowski comparison samples
# Append "Leak" suffix to the sample-type labels in column E for the
# leak-test rows (2-16) and the corresponding AU- rows (26-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..16) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $old = $cell.Value2
    $cell.Value = "$($old)Leak"
}

foreach ($r in 26..44) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $old = $cell.Value2
    $cell.Value = "$($old)Leak"
}

# Leave the selection on the last edited cell (E43), matching the
# recorded cursor position in the saved workbook.
$ws.Range("E43").Select()
